# Automatische test-sync: 2025-08-14 20:59:50
# Adds a new "CE-certificaten verzoek" log row (row 18) to the "Logs" sheet
# and bumps the "Intern verzoek / Actie voor medewerker" tally on "Dashboard".

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$row = 18
$logs.Cells.Item($row, 1).Value = "CE-certificaten verzoek"
$logs.Cells.Item($row, 2).Value = "inkoop@testbedrijf123.nl"
$logs.Cells.Item($row, 3).Value = "Kun je mij de CE-certificaten van de EcoPro-700 sturen?"
$logs.Cells.Item($row, 4).Value = "Intern verzoek / Actie voor medewerker"
$logs.Cells.Item($row, 5).Value = "Bedankt, we hebben dit doorgestuurd naar kwaliteit@testbedrijf123.nl."
$logs.Cells.Item($row, 6).Value = "2025-08-14 20:59:16"
$logs.Cells.Item($row, 7).Value = "Nee"
$logs.Cells.Item($row, 8).Value = "Ja"
$logs.Cells.Item($row, 9).Value = "Nee"
$logs.Cells.Item($row, 10).Value = "Nee"

$dashboard.Range("B2").Value = 12

# Extend the existing conditional-formatting ranges so the new row is covered too.
$logs.Range("D2:D17").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D18"))
$logs.Range("G2:G17").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G18"))
$logs.Range("H2:H17").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H18"))
$logs.Range("I2:I17").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I18"))
$logs.Range("J2:J17").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J18"))
